$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.942832
$ws.Range("H2").Value = 122.828496
$ws.Range("I2").Value = 0.2583000005785167
$ws.Range("J2").Value = 0.2583000005785167
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 967.7164930523253
$ws.Range("R2").Value = 8709.448437470928
$ws.Range("S2").Value = 0.01763632260932571
$ws.Range("T2").Value = 0.01763632260932571

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.942832
$ws.Range("H3").Value = 122.828496
$ws.Range("I3").Value = 0.2583000005785167
$ws.Range("J3").Value = 0.2583000005785167
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 7422.460149910924
$ws.Range("R3").Value = 66802.14134919831
$ws.Range("S3").Value = 0.135271954853016
$ws.Range("T3").Value = 0.135271954853016

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.942832
$ws.Range("H4").Value = 122.828496
$ws.Range("I4").Value = 0.2583000005785167
$ws.Range("J4").Value = 0.2583000005785167
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 4548.771003633883
$ws.Range("R4").Value = 40938.93903270494
$ws.Range("S4").Value = 0.08289989213989318
$ws.Range("T4").Value = 0.08289989213989318

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.942832
$ws.Range("H5").Value = 122.828496
$ws.Range("I5").Value = 0.2583000005785167
$ws.Range("J5").Value = 0.2583000005785167
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 1234.141395394059
$ws.Range("R5").Value = 11107.27255854653
$ws.Range("S5").Value = 0.02249183097628173
$ws.Range("T5").Value = 0.02249183097628173

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 79.68771233333334
$ws.Range("H6").Value = 239.063137
$ws.Range("I6").Value = 0.5027335710876245
$ws.Range("J6").Value = 0.5027335710876245
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 1883.482645230205
$ws.Range("R6").Value = 16951.34380707184
$ws.Range("S6").Value = 0.03432586692366102
$ws.Range("T6").Value = 0.03432586692366102

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 79.68771233333334
$ws.Range("H7").Value = 239.063137
$ws.Range("I7").Value = 0.5027335710876245
$ws.Range("J7").Value = 0.5027335710876245
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 14446.45717794343
$ws.Range("R7").Value = 130018.1146014908
$ws.Range("S7").Value = 0.2632820471503973
$ws.Range("T7").Value = 0.2632820471503972

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 79.68771233333334
$ws.Range("H8").Value = 239.063137
$ws.Range("I8").Value = 0.5027335710876245
$ws.Range("J8").Value = 0.5027335710876245
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 8853.348376286838
$ws.Range("R8").Value = 79680.13538658153
$ws.Range("S8").Value = 0.1613494337008287
$ws.Range("T8").Value = 0.1613494337008287

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 79.68771233333334
$ws.Range("H9").Value = 239.063137
$ws.Range("I9").Value = 0.5027335710876245
$ws.Range("J9").Value = 0.5027335710876245
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 2402.029847247019
$ws.Range("R9").Value = 21618.26862522317
$ws.Range("S9").Value = 0.04377622331273749
$ws.Range("T9").Value = 0.04377622331273749

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.53186833333334
$ws.Range("H10").Value = 112.595605
$ws.Range("I10").Value = 0.2367809244903433
$ws.Range("J10").Value = 0.2367809244903433
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 887.0956459786405
$ws.Range("R10").Value = 7983.860813807765
$ws.Range("S10").Value = 0.01616703353733329
$ws.Range("T10").Value = 0.01616703353733328

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 37.53186833333334
$ws.Range("H11").Value = 112.595605
$ws.Range("I11").Value = 0.2367809244903433
$ws.Range("J11").Value = 0.2367809244903433
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 6804.092033884475
$ws.Range("R11").Value = 61236.82830496028
$ws.Range("S11").Value = 0.1240023943320777
$ws.Range("T11").Value = 0.1240023943320776

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 37.53186833333334
$ws.Range("H12").Value = 112.595605
$ws.Range("I12").Value = 0.2367809244903433
$ws.Range("J12").Value = 0.2367809244903433
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 4169.811076744066
$ws.Range("R12").Value = 37528.29969069659
$ws.Range("S12").Value = 0.07599346905563362
$ws.Range("T12").Value = 0.0759934690556336

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 37.53186833333334
$ws.Range("H13").Value = 112.595605
$ws.Range("I13").Value = 0.2367809244903433
$ws.Range("J13").Value = 0.2367809244903433
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 1131.324583425155
$ws.Range("R13").Value = 10181.92125082639
$ws.Range("S13").Value = 0.02061802756529871
$ws.Range("T13").Value = 0.02061802756529871

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.3464216666666666
$ws.Range("H14").Value = 1.039265
$ws.Range("I14").Value = 0.002185503843515531
$ws.Range("J14").Value = 0.002185503843515531
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 8.187952420682777
$ws.Range("R14").Value = 73.69157178614499
$ws.Range("S14").Value = 0.000149222805891728
$ws.Range("T14").Value = 0.000149222805891728

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.3464216666666666
$ws.Range("H15").Value = 1.039265
$ws.Range("I15").Value = 0.002185503843515531
$ws.Range("J15").Value = 0.002185503843515531
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 62.80222667301222
$ws.Range("R15").Value = 565.2200400571099
$ws.Range("S15").Value = 0.001144550432013103
$ws.Range("T15").Value = 0.001144550432013103

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.3464216666666666
$ws.Range("H16").Value = 1.039265
$ws.Range("I16").Value = 0.002185503843515531
$ws.Range("J16").Value = 0.002185503843515531
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 38.48763642837055
$ws.Range("R16").Value = 346.3887278553349
$ws.Range("S16").Value = 0.0007014248257567696
$ws.Range("T16").Value = 0.0007014248257567696

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.3464216666666666
$ws.Range("H17").Value = 1.039265
$ws.Range("I17").Value = 0.002185503843515531
$ws.Range("J17").Value = 0.002185503843515531
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 10.44220192425222
$ws.Range("R17").Value = 93.97981731826998
$ws.Range("S17").Value = 0.0001903057798539309
$ws.Range("T17").Value = 0.0001903057798539309

